$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("latest")

$data = @{
    2  = @(-0.2453501895176297, 1.885822625406119, 16.01396840704841, 4.0017456699606, 4.083986181670692, 23)
    3  = @(-0.008306120572730035, 1.678411150852771, 10.3892641846745, 3.22323815202577, 3.299078473400757, 22)
    4  = @(-0.6206486848402303, 1.081319350267402, 4.95937692841922, 2.226965857039398, 2.191547568579468, 21)
    5  = @(-0.007917461514596722, 0.7359397285945519, 2.052934180666193, 1.432806400274019, 1.470005905634166, 20)
    6  = @(-0.01396789604063663, 0.6794082136444377, 0.985853761070935, 0.992901687515403, 1.020008565325462, 19)
    7  = @(0.006227155691902452, 0.5859396340579296, 0.8714568099673092, 0.9335185107791432, 0.9605612509173802, 18)
    8  = @(0.1323879842408515, 0.5124844762799831, 0.586790449754057, 0.7660224864545799, 0.7777164414773035, 17)
    9  = @(0.1676475329133778, 0.4867046202604076, 0.3917272528104276, 0.6258811810642877, 0.6227864759983028, 16)
    10 = @(0.199895990469669, 0.4747094579524145, 0.3665317847062111, 0.6054186854617316, 0.5915234970278875, 15)
    11 = @(0.263452728002003, 0.3988499978544717, 0.2138270447762498, 0.4624143648030949, 0.3943715417630878, 14)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
